$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: list of (address, new value) pairs reflecting the refreshed
# coinranking snapshot captured by the scheduled GitHub Actions job.
$updates = @(
    @('D2', '296.80'),
    @('E2', '2.16%'),
    @('G2', '23'),
    @('D3', '41.95'),
    @('E3', '3.61%'),
    @('G3', '23'),
    @('D4', '5.039'),
    @('E4', '0.39%'),
    @('G4', '23'),
    @('D5', '0.07571'),
    @('E5', '3.33%'),
    @('G5', '23'),
    @('B6', 'GateToken'),
    @('C6', 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'),
    @('D6', '4.398'),
    @('E6', '2.74%'),
    @('G6', '23'),
    @('B7', 'FTXToken'),
    @('C7', 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'),
    @('D7', '1.602'),
    @('E7', '2.69%'),
    @('G7', '23'),
    @('B8', 'MXToken'),
    @('C8', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'),
    @('D8', '0.9288'),
    @('E8', '0.54%'),
    @('G8', '23'),
    @('B9', 'BTSEToken'),
    @('C9', 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'),
    @('D9', '2.408'),
    @('E9', '3.44%'),
    @('G9', '23'),
    @('B10', 'LiechtensteinCryptoassetsExchange'),
    @('C10', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'),
    @('D10', '0.1209'),
    @('E10', '5.01%'),
    @('G10', '23'),
    @('B11', 'WazirX'),
    @('C11', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'),
    @('D11', '0.1836'),
    @('E11', '5.72%'),
    @('G11', '23'),
    @('B12', 'MandalaExchangeToken'),
    @('C12', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'),
    @('D12', '0.09005'),
    @('E12', '4.89%'),
    @('G12', '23'),
    @('B13', 'BitrueCoin'),
    @('C13', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'),
    @('D13', '0.04019'),
    @('E13', '-3.72%'),
    @('G13', '23'),
    @('B14', 'BitMartToken'),
    @('C14', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'),
    @('D14', '0.1052'),
    @('E14', '-0.20%'),
    @('G14', '23'),
    @('B15', 'BitForexToken'),
    @('C15', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'),
    @('D15', '0.001279'),
    @('E15', '0.10%'),
    @('G15', '23'),
    @('B16', 'TigerCash'),
    @('C16', 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'),
    @('D16', '0.005822'),
    @('E16', '-1.60%'),
    @('G16', '23'),
    @('B17', 'LEO'),
    @('C17', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'),
    @('D17', '3.363'),
    @('E17', '-1.47%'),
    @('G17', '23'),
    @('D18', '0.3321'),
    @('E18', '1.17%'),
    @('G18', '23'),
    @('D19', '7.906'),
    @('E19', '1.13%'),
    @('G19', '23'),
    @('D20', '0.1420'),
    @('E20', '2.91%'),
    @('G20', '23'),
    @('D21', '0.3002'),
    @('E21', '4.18%'),
    @('G21', '23'),
    @('D22', '0.04062'),
    @('E22', '5.28%'),
    @('G22', '23'),
    @('D23', '0.001265'),
    @('E23', '0.39%'),
    @('G23', '23'),
    @('D24', '0.003951'),
    @('E24', '4.07%'),
    @('G24', '23'),
    @('D25', '0.0001229'),
    @('E25', '-4.07%'),
    @('G25', '23'),
    @('E26', '0.03%'),
    @('G26', '23'),
    @('G27', '23'),
    @('G28', '23'),
    @('G29', '23'),
    @('G30', '23'),
    @('G31', '23'),
    @('G32', '23'),
    @('G33', '23'),
    @('G34', '23'),
    @('G35', '23'),
    @('G36', '23'),
    @('G37', '23'),
    @('D38', '0.02413'),
    @('E38', '3.97%'),
    @('G38', '23'),
    @('D39', '0.05215'),
    @('G39', '23'),
    @('D40', '0.005946'),
    @('E40', '-9.07%'),
    @('G40', '23'),
    @('D41', '0.007770'),
    @('E41', '1.01%'),
    @('G41', '23'),
    @('D42', '0.1330'),
    @('E42', '4.09%'),
    @('G42', '23'),
    @('D43', '0.007537'),
    @('E43', '2.47%'),
    @('G43', '23'),
    @('D44', '0.007845'),
    @('E44', '10.79%'),
    @('G44', '23'),
    @('D45', '0.2972'),
    @('E45', '-5.80%'),
    @('G45', '23'),
    @('D46', '0.00006788'),
    @('E46', '6.14%'),
    @('G46', '23'),
    @('E47', '0.05%'),
    @('G47', '23'),
    @('D48', '0.04551'),
    @('E48', '183.38%'),
    @('G48', '23'),
    @('D49', '0.004203'),
    @('E49', '0.03%'),
    @('G49', '23'),
    @('D50', '0.00002102'),
    @('E50', '0.05%'),
    @('G50', '23'),
    @('D51', '0.0002002'),
    @('E51', '0.05%'),
    @('G51', '23')
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

